$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) values that look numeric are stored as text, like the source data.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Rows 2-45: in-place Price / Volume(1h) updates ---
$ws.Cells.Item(2, 4).Value = '26.365.86'
$ws.Cells.Item(2, 5).Value = '  +1.32%  '
$ws.Cells.Item(3, 4).Value = '1.622.78'
$ws.Cells.Item(3, 5).Value = '  +1.65%  '
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$ws.Cells.Item(5, 4).Value = '212.28'
$ws.Cells.Item(5, 5).Value = '  +0.23%  '
$ws.Cells.Item(6, 5).Value = '  -0.09%  '
$ws.Cells.Item(7, 5).Value = '  +1.32%  '
$ws.Cells.Item(8, 5).Value = '  +0.08%  '
$ws.Cells.Item(9, 4).Value = '0.0617'
$ws.Cells.Item(9, 5).Value = '  +0.44%  '
$ws.Cells.Item(10, 4).Value = '18.89'
$ws.Cells.Item(10, 5).Value = '  +3.64%  '
$ws.Cells.Item(11, 5).Value = '  +0.81%  '
$ws.Cells.Item(12, 4).Value = '1.848.65'
$ws.Cells.Item(13, 4).Value = '1.624.91'
$ws.Cells.Item(13, 5).Value = '  +1.86%  '
$ws.Cells.Item(14, 5).Value = '  +0.26%  '
$ws.Cells.Item(15, 5).Value = '  +0.73%  '
$ws.Cells.Item(16, 4).Value = '26.365.96'
$ws.Cells.Item(16, 5).Value = '  +1.35%  '
$ws.Cells.Item(17, 4).Value = '62.52'
$ws.Cells.Item(17, 5).Value = '  +2.90%  '
$ws.Cells.Item(18, 5).Value = '  -0.17%  '
$ws.Cells.Item(19, 5).Value = '  -0.07%  '
$ws.Cells.Item(20, 4).Value = '202.99'
$ws.Cells.Item(20, 5).Value = '  -0.38%  '
$ws.Cells.Item(21, 5).Value = '  +0.27%  '
$ws.Cells.Item(22, 5).Value = '  +0.60%  '
$ws.Cells.Item(23, 5).Value = '  +0.00%  '
$ws.Cells.Item(24, 5).Value = '  -2.53%  '
$ws.Cells.Item(25, 4).Value = '144.99'
$ws.Cells.Item(25, 5).Value = '  +0.72%  '
$ws.Cells.Item(26, 5).Value = '  -0.01%  '
$ws.Cells.Item(27, 5).Value = '  -0.76%  '
$ws.Cells.Item(28, 5).Value = '  +0.13%  '
$ws.Cells.Item(29, 5).Value = '  +0.91%  '
$ws.Cells.Item(30, 4).Value = '0.0522'
$ws.Cells.Item(30, 5).Value = '  +9.68%  '
$ws.Cells.Item(31, 5).Value = '  -0.16%  '
$ws.Cells.Item(32, 5).Value = '  +1.84%  '
$ws.Cells.Item(33, 5).Value = '  +0.42%  '
$ws.Cells.Item(34, 4).Value = '1.49'
$ws.Cells.Item(34, 5).Value = '  +0.69%  '
$ws.Cells.Item(35, 5).Value = '  +2.15%  '
$ws.Cells.Item(36, 4).Value = '1.178.16'
$ws.Cells.Item(36, 5).Value = '  +4.25%  '
$ws.Cells.Item(37, 5).Value = '  -0.12%  '
$ws.Cells.Item(38, 4).Value = '0.809'
$ws.Cells.Item(38, 5).Value = '  +1.42%  '
$ws.Cells.Item(39, 5).Value = '  -0.05%  '
$ws.Cells.Item(40, 4).Value = '2.33'
$ws.Cells.Item(40, 5).Value = '  -0.04%  '
$ws.Cells.Item(41, 4).Value = '0.499'
$ws.Cells.Item(41, 5).Value = '  +1.19%  '
$ws.Cells.Item(42, 5).Value = '  +4.41%  '
$ws.Cells.Item(43, 4).Value = '0.785'
$ws.Cells.Item(43, 5).Value = '  +0.66%  '
$ws.Cells.Item(44, 4).Value = '1.759.90'
$ws.Cells.Item(44, 5).Value = '  +1.75%  '
$ws.Cells.Item(45, 4).Value = '92.62'
$ws.Cells.Item(45, 5).Value = '  +0.34%  '

# --- New coin (BabyDogeCoin) inserted at row 46; existing rows 46-50 shift down to 47-51; ---
# --- the former last row (EnergySwap, row 51) drops off the bottom of the 51-row table. ---
$ws.Rows.Item(46).Insert()

# Restore the row-header (column A) formatting that Insert() did not carry over cleanly.
$ws.Cells.Item(45, 1).Copy()
$ws.Cells.Item(46, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(46, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(46, 4).Value = '0.0₆0104'
$ws.Cells.Item(46, 5).Value = '  +9.94%  '

# --- Rows 47-51 now hold the coins that were at 46-50 before the insert; refresh their rank (A),
# --- plus Price / Volume(1h). ---
$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 5).Value = '  +2.29%  '
$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 4).Value = '53.83'
$ws.Cells.Item(48, 5).Value = '  -0.32%  '
$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(49, 5).Value = '  +0.58%  '
$ws.Cells.Item(50, 1).Value = 48
$ws.Cells.Item(50, 5).Value = '  +0.79%  '
$ws.Cells.Item(51, 1).Value = 49
$ws.Cells.Item(51, 5).Value = '  -0.40%  '

# Drop the row that fell off the bottom of the table (old EnergySwap row, now duplicated at 52).
$ws.Rows.Item(52).Delete()
